$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.816.19'
$ws.Range('E2').Value = '  +1.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.662.21'
$ws.Range('E3').Value = '  +3.09%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.49'
$ws.Range('E5').Value = '  +2.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.85'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  -0.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.64'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.66'
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.138.00'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.745.28'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.662.62'
$ws.Range('E17').Value = '  +3.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.43'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '343.68'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.38'
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.79'
$ws.Range('E21').Value = '  +1.75%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.15'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.63'
$ws.Range('E24').Value = '  +13.73%  '
$ws.Range('E25').Value = '  +5.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '577.73'
$ws.Range('E26').Value = '  +24.57%  '
$ws.Range('E27').Value = '  -1.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.06'
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.53'
$ws.Range('E29').Value = '  +2.59%  '
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('E31').Value = '  +3.95%  '
$ws.Range('E32').Value = '  +12.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0₃0819'
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '175.46'
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.79'
$ws.Range('E37').Value = '  +5.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.21'
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('E39').Value = '  +3.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '171.40'
$ws.Range('E40').Value = '  +8.44%  '
$ws.Range('E42').Value = '  +2.82%  '
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('E44').Value = '  +3.20%  '
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0553'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0240'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0964'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.72'
$ws.Range('E49').Value = '  +1.57%  '
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0216'
$ws.Range('E51').Value = '  +11.34%  '
